$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.714.93"
$c.ClearFormats()
$ws.Range("E2").Value = "  +0.75%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.845.91"
$c.ClearFormats()
$ws.Range("E3").Value = "  +0.17%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.010"
$c.ClearFormats()
$ws.Range("E4").Value = "  +0.35%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "335.64"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("E6").Value = "  +0.25%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4651"
$c.ClearFormats()
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3854"
$c.ClearFormats()
$ws.Range("E8").Value = "  -0.27%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "46.77"
$c.ClearFormats()
$ws.Range("E9").Value = "  +1.78%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07893"
$c.ClearFormats()
$ws.Range("E10").Value = "  -0.22%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.9662"
$c.ClearFormats()
$ws.Range("E11").Value = "  -3.06%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "21.22"
$c.ClearFormats()
$ws.Range("E12").Value = "  -1.32%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.861.85"
$c.ClearFormats()
$ws.Range("E13").Value = "  +0.71%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.882"
$c.ClearFormats()
$ws.Range("E14").Value = "  -1.37%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.132"
$c.ClearFormats()
$ws.Range("E15").Value = "  -0.01%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.009"
$c.ClearFormats()
$ws.Range("E16").Value = "  +0.13%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "90.93"
$c.ClearFormats()
$ws.Range("E17").Value = "  +2.99%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.06619"
$c.ClearFormats()
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("E19").Value = "  -0.49%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "17.26"
$c.ClearFormats()
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("E21").Value = "  +0.16%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "27.728.26"
$c.ClearFormats()
$ws.Range("E22").Value = "  +0.78%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.346"
$c.ClearFormats()
$ws.Range("E23").Value = "  -0.85%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "10.80"
$c.ClearFormats()
$ws.Range("E24").Value = "  -0.75%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.303"
$c.ClearFormats()
$ws.Range("E25").Value = "  -0.25%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.082.31"
$c.ClearFormats()
$ws.Range("E26").Value = "  +0.73%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "158.91"
$c.ClearFormats()
$ws.Range("E27").Value = "  -0.01%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "19.43"
$c.ClearFormats()
$ws.Range("E28").Value = "  -0.31%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.062"
$c.ClearFormats()
$ws.Range("E29").Value = "  -2.27%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "5.386"
$c.ClearFormats()
$ws.Range("E30").Value = "  -0.29%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "118.65"
$c.ClearFormats()
$ws.Range("E31").Value = "  -1.19%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.09436"
$c.ClearFormats()
$ws.Range("E32").Value = "  +0.46%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.9391"
$c.ClearFormats()
$ws.Range("E33").Value = "  -3.58%  "
$ws.Range("E34").Value = "  +0.29%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.254"
$c.ClearFormats()
$ws.Range("E35").Value = "  -0.76%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.324"
$c.ClearFormats()
$ws.Range("E36").Value = "  -0.77%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.06022"
$c.ClearFormats()
$ws.Range("E37").Value = "  +0.16%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02208"
$c.ClearFormats()
$ws.Range("E38").Value = "  -0.71%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "8.216"
$c.ClearFormats()
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("E40").Value = "  +0.08%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.150"
$c.ClearFormats()
$ws.Range("E41").Value = "  -2.52%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.5794"
$c.ClearFormats()
$ws.Range("E42").Value = "  -1.77%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.1842"
$c.ClearFormats()
$ws.Range("E43").Value = "  -0.94%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "10.04"
$c.ClearFormats()
$ws.Range("E44").Value = "  -2.95%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.300"
$c.ClearFormats()
$ws.Range("E45").Value = "  +4.96%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "11.98"
$c.ClearFormats()
$ws.Range("E46").Value = "  -0.86%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.5440"
$c.ClearFormats()
$ws.Range("E47").Value = "  -2.55%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.933"
$c.ClearFormats()
$ws.Range("E48").Value = "  +1.30%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.06836"
$c.ClearFormats()
$ws.Range("E49").Value = "  +2.19%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "110.66"
$c.ClearFormats()
$ws.Range("E50").Value = "  +0.60%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.009"
$c.ClearFormats()
$ws.Range("E51").Value = "  -32.21%  "
